$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the paragraph "What is the main Diagram of UML classes modeling everything?"
$findRange = $d.Content
[void]$findRange.Find.Execute("What is the main Diagram of UML classes modeling everything?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$diagramStart = $findRange.Paragraphs.Item(1).Range.Start

$diagramIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $diagramStart) {
        $diagramIndex = $i
    }
}

# 1. Insert a new paragraph BEFORE it with the new "structures" question.
$beforeRange = $d.Paragraphs.Item($diagramIndex).Range
$beforeRange.Collapse(1)
$beforeRange.InsertParagraphBefore()
$diagramIndex = $diagramIndex + 1
$d.Paragraphs.Item($diagramIndex - 1).Range.Text = "What are the structures being used, and how do they relate to each other?"

# 2. After the "Diagram" paragraph, insert three more paragraphs:
#    - an empty paragraph (no run at all)
#    - a paragraph explaining the project's functionalities
#    - a paragraph (two runs) about switching algorithms/sources
$afterRange = $d.Paragraphs.Item($diagramIndex).Range
$afterRange.Collapse(0)
$afterRange.InsertParagraphAfter()
$afterRange.InsertParagraphAfter()
$afterRange.InsertParagraphAfter()

$emptyIndex = $diagramIndex + 1
$funcIndex = $diagramIndex + 2
$switchIndex = $diagramIndex + 3

[void]$d.Paragraphs.Item($emptyIndex).Range.InsertXML("<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr></w:p>")

$d.Paragraphs.Item($funcIndex).Range.Text = "I have to explain here also the different functionalities of what the project does: how to process a template or trace form a repository (and which classes to extend), how to change a domain, how to save the files in a format, how produce the script to actually run SUBDUE or any other graph, how to create statistics, how to validate the fragments found, how to generate the fragment catalog in wffd, etc."

$switchXml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t xml:space=`"preserve`">How to switch from different </w:t></w:r><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>algorithms and how to attack different sources (which is in fact easy).</w:t></w:r></w:p>"
[void]$d.Paragraphs.Item($switchIndex).Range.InsertXML($switchXml)
